$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary cell ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.08 = 33441.84 pesos`n✅ 33441.84 pesos = 8.03 = 957.02 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 123.797
$wsTasas.Range("O10").Value = 4140
$wsTasas.Range("N12").Value = 4166
$wsTasas.Range("O12").Value = 119.22
